$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.120.25"
$ws.Cells.Item(2, 5).Value = "  -2.22%  "

$ws.Cells.Item(3, 4).Value = "1.558.46"
$ws.Cells.Item(3, 5).Value = "  -2.32%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "206.33"
$ws.Cells.Item(5, 5).Value = "  -1.02%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.488"
$ws.Cells.Item(6, 5).Value = "  -2.78%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "22.07"
$ws.Cells.Item(8, 5).Value = "  -1.39%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.247"
$ws.Cells.Item(9, 5).Value = "  -2.36%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0591"
$ws.Cells.Item(10, 5).Value = "  -0.46%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0865"
$ws.Cells.Item(11, 5).Value = "  -0.37%  "

$ws.Cells.Item(12, 4).Value = "1.781.35"
$ws.Cells.Item(12, 5).Value = "  -2.19%  "

$ws.Cells.Item(13, 4).Value = "1.565.41"
$ws.Cells.Item(13, 5).Value = "  -2.40%  "

$ws.Cells.Item(14, 5).Value = "  -2.87%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.515"
$ws.Cells.Item(15, 5).Value = "  -3.42%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "62.89"
$ws.Cells.Item(16, 5).Value = "  -1.03%  "

$ws.Cells.Item(17, 4).Value = "27.128.25"
$ws.Cells.Item(17, 5).Value = "  -2.19%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "214.81"
$ws.Cells.Item(18, 5).Value = "  -2.40%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0685"
$ws.Cells.Item(19, 5).Value = "  -1.84%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.21"
$ws.Cells.Item(20, 5).Value = "  -2.32%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.11"
$ws.Cells.Item(22, 5).Value = "  -1.26%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "9.31"
$ws.Cells.Item(23, 5).Value = "  -4.41%  "

$ws.Cells.Item(24, 5).Value = "  +0.37%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "151.55"
$ws.Cells.Item(25, 5).Value = "  -1.57%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.58"
$ws.Cells.Item(26, 5).Value = "  -3.12%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "14.91"
$ws.Cells.Item(27, 5).Value = "  -1.75%  "

$ws.Cells.Item(28, 5).Value = "  +0.05%  "

$ws.Cells.Item(29, 5).Value = "  -1.66%  "

$ws.Cells.Item(31, 5).Value = "  -2.48%  "

$ws.Cells.Item(32, 5).Value = "  -2.49%  "

$ws.Cells.Item(33, 4).Value = "1.380.68"
$ws.Cells.Item(33, 5).Value = "  +0.29%  "

$ws.Cells.Item(34, 5).Value = "  -1.21%  "

$ws.Cells.Item(35, 5).Value = "  -0.47%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.944"
$ws.Cells.Item(36, 5).Value = "  -2.84%  "

$ws.Cells.Item(37, 5).Value = "  -1.74%  "

$ws.Cells.Item(38, 5).Value = "  -1.92%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.809"
$ws.Cells.Item(39, 5).Value = "  -2.43%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.512"
$ws.Cells.Item(40, 5).Value = "  -4.75%  "

$ws.Cells.Item(41, 5).Value = "  +0.05%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.988"
$ws.Cells.Item(42, 5).Value = "  +1.52%  "

$ws.Cells.Item(43, 5).Value = "  +3.03%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "63.26"
$ws.Cells.Item(44, 5).Value = "  -2.11%  "

$ws.Cells.Item(45, 5).Value = "  -0.45%  "

$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.21"
$ws.Cells.Item(46, 5).Value = "  +0.12%  "

$ws.Cells.Item(47, 4).Value = "1.694.25"
$ws.Cells.Item(47, 5).Value = "  -2.11%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "85.31"
$ws.Cells.Item(48, 5).Value = "  -1.97%  "

$ws.Cells.Item(49, 5).Value = "  -3.25%  "

$ws.Cells.Item(50, 5).Value = "  -1.10%  "

$ws.Cells.Item(51, 5).Value = "  +0.02%  "
